# Weekly update: a new price record is inserted at row 35 (pushing the
# existing rows 35-66 down to 36-67), matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 35; Excel shifts rows 35..66 down to 36..67 and
# carries the row format (so D35 keeps its date number format).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = 44893
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = 100112026
$ws.Cells.Item(35, 7).Value = "Haba"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 120
$ws.Cells.Item(35, 11).Value = 9000
$ws.Cells.Item(35, 12).Value = 9500
$ws.Cells.Item(35, 13).Value = 9250
$ws.Cells.Item(35, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(35, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(35, 16).Value = 370
$ws.Cells.Item(35, 17).Value = 25
$ws.Cells.Item(35, 18).Value = "Hortaliza"
